# Apply updated "want to go" counts (F) and min-price (G) values
# from the site re-scrape, per sheet/row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 107
$ws.Range("F6").Value = 1053
$ws.Range("F9").Value = 11828
$ws.Range("F13").Value = 2125
$ws.Range("F18").Value = 1207
$ws.Range("F19").Value = 179
$ws.Range("F20").Value = 256
$ws.Range("F23").Value = 279
$ws.Range("F24").Value = 2907
$ws.Range("F25").Value = 738
$ws.Range("F26").Value = 3712
$ws.Range("F27").Value = 3712
$ws.Range("F28").Value = 1080
$ws.Range("F29").Value = 824
$ws.Range("G29").Value = 58.5
$ws.Range("F32").Value = 21
$ws.Range("F33").Value = 1004
$ws.Range("F34").Value = 41
$ws.Range("F35").Value = 83
$ws.Range("F36").Value = 260
$ws.Range("F40").Value = 3930
$ws.Range("F41").Value = 4461
$ws.Range("F42").Value = 5498
$ws.Range("F44").Value = 119
$ws.Range("F45").Value = 164
$ws.Range("F47").Value = 70

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4162
$ws.Range("F12").Value = 781
$ws.Range("F20").Value = 80

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 756

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 756
$ws.Range("F8").Value = 1053
$ws.Range("F11").Value = 11828
$ws.Range("F14").Value = 2125
$ws.Range("F18").Value = 1207
$ws.Range("F19").Value = 179
$ws.Range("F20").Value = 256
$ws.Range("F21").Value = 4162
$ws.Range("F23").Value = 279
$ws.Range("F24").Value = 738
$ws.Range("F25").Value = 3712
$ws.Range("F26").Value = 1080
$ws.Range("F29").Value = 824
$ws.Range("G29").Value = 58.5
$ws.Range("F30").Value = 1004
$ws.Range("F31").Value = 41
$ws.Range("F32").Value = 83
$ws.Range("F33").Value = 260
$ws.Range("F36").Value = 4461
$ws.Range("F38").Value = 119
$ws.Range("F39").Value = 164
$ws.Range("F44").Value = 70
$ws.Range("F48").Value = 80
